# Daily attendance processing - normalize the "Recorded By" (column G) lists.
# Rule: for each "Recorded By" cell that lists multiple comma-separated
# recorders, move the exact ("System", case-sensitive) entry to the front
# (preserving the relative order of the remaining entries). If no exact
# "System" entry is present, sort the entries alphabetically instead.
# Single-entry cells are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $orig = $cell.Text

    if ([string]::IsNullOrEmpty($orig)) {
        continue
    }

    $parts = $orig.Split(",")
    if ($parts.Count -le 1) {
        continue
    }

    $trimmedParts = @()
    foreach ($p in $parts) {
        $trimmedParts += $p.Trim()
    }

    # Locate the first *exact-case* "System" entry.
    $sysIdx = -1
    for ($i = 0; $i -lt $trimmedParts.Count; $i++) {
        if ($trimmedParts[$i].Equals("System")) {
            $sysIdx = $i
            break
        }
    }

    if ($sysIdx -ge 0) {
        $rest = @()
        for ($i = 0; $i -lt $trimmedParts.Count; $i++) {
            if ($i -ne $sysIdx) {
                $rest += $trimmedParts[$i]
            }
        }
        $newParts = @("System") + $rest
    } else {
        $newParts = $trimmedParts | Sort-Object
    }

    $newVal = $newParts -join ", "

    if (-not $newVal.Equals($orig)) {
        $cell.Value = $newVal
    }
}
